## maintenance + adminPage + error_log
## Refresh the "staffEmail" column (E) of the Banquet test-data sheet with
## a new batch of sample addresses (some rows now share the regenerated
## obfuscated mailbox names used by the admin-page / error-log fixtures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$emails = @{
    2  = "emily.johnson@email.com"
    3  = "jane.brown@example.com"
    4  = "risxg.wfyfd@example.com"
    5  = "li.johnson@email.com"
    6  = "li.chen@example.com"
    7  = "li.anderson@email.com"
    8  = "natalie.smith@testmail.com"
    9  = "srriu.wzhbr@example.com"
    10 = "john.doe@testmail.com"
    11 = "tzewd.mabrc@example.com"
    12 = "neerd.adsfr@example.com"
    13 = "jaqwt.uyzak@example.com"
    14 = "emily.johnson@email.com"
    15 = "jane.williams@testmail.com"
    16 = "cicaj.lzfgk@example.com"
    17 = "zmmrk.iskdu@example.com"
    18 = "grace.anderson@email.com"
    19 = "li.doe@testmail.com"
    20 = "chris.doe@email.com"
    21 = "mark.williams@testmail.com"
    22 = "john.williams@testmail.com"
    23 = "john.williams@example.com"
    24 = "li.lee@example.com"
    25 = "jane.lee@testmail.com"
    26 = "tom.patel@testmail.com"
    27 = "mark.lee@email.com"
    28 = "mark.smith@testmail.com"
    29 = "grace.chen@testmail.com"
    30 = "amit.smith@example.com"
    31 = "grace.smith@example.com"
    32 = "jane.clark@email.com"
    33 = "sophia.patel@email.com"
    34 = "li.patel@testmail.com"
    35 = "amit.doe@testmail.com"
    36 = "logkk.visap@example.com"
    37 = "kezfr.pwqqq@example.com"
    38 = "dczaw.xahwj@example.com"
    39 = "anna.johnson@example.com"
    40 = "jane.chen@testmail.com"
    41 = "siquo.jqrgv@example.com"
    42 = "grace.williams@example.com"
    43 = "emily.patel@example.com"
    44 = "amit.clark@example.com"
    45 = "jjkua.jedeq@example.com"
    46 = "mark.williams@example.com"
    47 = "qcpdl.hmqvc@example.com"
    48 = "zabcp.xhyeg@example.com"
    49 = "natalie.clark@email.com"
    50 = "amit.chen@example.com"
    51 = "vsuwx.tkmnt@example.com"
}

foreach ($row in $emails.Keys) {
    $ws.Range("E$row").Value = $emails[$row]
}

# Match the author's final on-screen selection.
$ws.Range("E51").Select()
